$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new label value (column A), rows 2-38.
# This re-orders the emotion/area labels that were previously scrambled
# (bug with mixing emotions in cache).
$values = @{
    2  = "51-1-1"
    3  = "47-4-2"
    4  = "50-1-2"
    5  = "45-5-2"
    6  = "49-2-1"
    7  = "41-3-1"
    8  = "50-2-1"
    9  = "45-5-1"
    10 = "50-3-2"
    11 = "37-1-1"
    12 = "49-3-2"
    13 = "40-2-2"
    14 = "49-2-2"
    15 = "38-1-1"
    16 = "49-1-1"
    17 = "40-5-1"
    18 = "50-3-1"
    19 = "56-2-2"
    20 = "42-3-2"
    21 = "47-1-1"
    22 = "51-2-2"
    23 = "46-3-1"
    24 = "46-3-2"
    25 = "49-1-2"
    26 = "50-4-1"
    27 = "66-3-1"
    28 = "30-4-3"
    29 = "52-2-1"
    30 = "37-5-1"
    31 = "46-4-1"
    32 = "52-3-1"
    33 = "38-2-1"
    34 = "46-1-2"
    35 = "48-4-1"
    36 = "52-2-2"
    37 = "38-4-1"
    38 = "40-4-2"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
